# "Fourth column of data" — fill in column G (rows 6-85) on the Time
# Analysis sheet with the (previously blank) data that was computed for
# the fourth magnitude bucket, and leave the selection where the author
# left it (cell G86, having scrolled down near the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = [ordered]@{
    6  = 5775
    7  = 5828
    8  = 3250
    9  = 3214
    10 = 3199
    11 = 22
    12 = 22
    13 = 19
    14 = 16
    15 = 17
    16 = 24234
    17 = 24313
    18 = 12675
    19 = 12459
    20 = 12523
    21 = 20
    22 = 22
    23 = 14
    24 = 12
    25 = 12
    26 = 1280
    27 = 1206
    28 = 2125
    29 = 2371
    30 = 2461
    31 = 28
    32 = 21
    33 = 19
    34 = 17
    35 = 18
    36 = 22
    37 = 15
    38 = 10
    39 = 10
    40 = 11
    41 = 234
    42 = 34
    43 = 37
    44 = 40
    45 = 41
    46 = 41
    47 = 18
    48 = 16
    49 = 10
    50 = 9
    51 = 6
    52 = 2
    53 = 1
    54 = 0
    55 = 1
    56 = 3
    57 = 2
    58 = 1
    59 = 1
    60 = 1
    61 = 3
    62 = 0
    63 = 1
    64 = 1
    65 = 0
    66 = 20
    67 = 13
    68 = 11
    69 = 10
    70 = 10
    71 = 8
    72 = 8
    73 = 2
    74 = 2
    75 = 2
    76 = 6
    77 = 6
    78 = 2
    79 = 1
    80 = 2
    81 = 6
    82 = 4
    83 = 2
    84 = 2
    85 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# Match the saved workbook's final view: scrolled near the bottom of the
# table with G86 as the active cell.
$ws.Range("G86").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 69
